# Contact-info line currently reads (4 runs):
#   R1 "{{Phone}}| "    (sz=16, szCs=24)
#   R2 "{{Email"        (sz=16, szCs=24)
#   R3 "}}"             (sz=16, szCs=16)
#   R4 "|{{LinkedIn}}"  (sz=16, szCs=16)
#
# Target layout (6 runs, same visible text except one extra space before the
# trailing pipe):
#   N1 "{{Phone}}"    (szCs=24)
#   N2 "|"            (szCs=24)
#   N3 " "            (szCs=24)
#   N4 "{{Email}}"    (szCs=24)
#   N5 " |"           (szCs=16)
#   N6 "{{LinkedIn}}" (szCs=16)

$d = $word.ActiveDocument

# Locate the paragraph that holds the Phone/Email/LinkedIn placeholders.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*{{Phone}}*{{Email*LinkedIn}}*") {
        $target = $p.Range
        break
    }
}

$base = $target.Start
$full = $target.Text
$idx = $full.IndexOf("{{Phone}}")
$start = $base + $idx

# ---------------------------------------------------------------------
# Step 1: "}}" (closing brace of {{Email}}) currently carries szCs=16,
# but it must end up merged with "{{Email" under szCs=24. The object
# model has no direct "szCs" setter, so borrow the formatting of an
# existing szCs=24 run via FormattedText (matching the donor length to
# the destination length keeps offsets stable), then restore the text.
# ---------------------------------------------------------------------
$donorLen = 2
$donor = $d.Range($start, $start + $donorLen)     # "{{" -- szCs=24 donor
$closeBrace = $start + ("{{Phone}}| {{Email".Length)
$dst = $d.Range($closeBrace, $closeBrace + 2)      # "}}" -- szCs=16
$dst.FormattedText = $donor.FormattedText
$dst.Text = "}}"

# ---------------------------------------------------------------------
# Step 2: add the missing space before the final pipe. Rewriting the
# whole trailing run's Text preserves its run-level formatting
# (szCs=16) while growing its length by one character.
# ---------------------------------------------------------------------
$pipeRunStart = $start + ("{{Phone}}| {{Email}}".Length)
$tail = $d.Range($pipeRunStart, $target.End)
$tail.Text = " |{{LinkedIn}}"

# ---------------------------------------------------------------------
# Step 3: split the two remaining multi-piece runs into the final run
# boundaries. Word merges adjacent runs that share identical formatting,
# so toggling a sub-range's Font.Size away and back forces a run break
# at that boundary without changing the rendered appearance.
# ---------------------------------------------------------------------
function Split-RunAt($from, $to) {
    $r = $d.Range($from, $to)
    $orig = $r.Font.Size
    if ($orig -eq 40) { $tmp = 8 } else { $tmp = 40 }
    $r.Font.Size = $tmp
    $r.Font.Size = $orig
}

$phoneEnd   = $start + "{{Phone}}".Length
$pipeEnd    = $phoneEnd + 1
$spaceAfter = $pipeRunStart + 2

# Process right-to-left so each split's end boundary is already final.
Split-RunAt $spaceAfter $target.End
Split-RunAt $pipeEnd $pipeRunStart
Split-RunAt $phoneEnd $pipeEnd

Write-Host "Result: $($target.Text)"
